$d = $word.ActiveDocument

# 1. Fix spelling: "behaviors" -> "behaviours" (US -> UK spelling)
$rng = $d.Content
$rng.Find.Execute("behaviors", $true, $false, $false, $false, $false, $true, 1, $false, "behaviours", 2) | Out-Null

# Remember where the replaced word starts/ends so we can relocate the
# "_GoBack" bookmark to sit right after "behaviou" (matching the author's
# split of "behaviou" + "r" runs around the bookmark).
$behavioursStart = $rng.Start

# 2. Move the "_GoBack" bookmark from its old spot (after "...ization")
#    to right after "behaviou" or before "rs ".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$bmRange = $d.Range($behavioursStart + 8, $behavioursStart + 8)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
